$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing data (D->E ... K->L)
$ws.Range("D1").EntireColumn.Insert()

# New column D inherits default style; copy number formats/styles from column E (the old column D)
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new (most-recent-period) figures in column D
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 212500
$ws.Range("D9").Value = 51700
$ws.Range("D10").Value = 160800
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 300
$ws.Range("D15").Value = 80000
$ws.Range("D17").Value = 157200
$ws.Range("D18").Value = 55300
$ws.Range("D20").Value = 17200
$ws.Range("D21").Value = 152500
$ws.Range("D22").Value = 25400
$ws.Range("D23").Value = 47100
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 47100
$ws.Range("D27").Value = 36100
$ws.Range("D28").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -17200
$ws.Range("D33").Value = 36100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 36100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 180600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 28100
$ws.Range("D48").Value = 2487300
$ws.Range("D49").Value = 60800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2787700
$ws.Range("D57").Value = 21100
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 39200
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 757400
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 911600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 159100
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -88300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1717000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 36100
$ws.Range("D83").Value = 80000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 102800
$ws.Range("D91").Value = -58000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -507200
$ws.Range("D96").Value = -65500
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 578200
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 173700

# Rows where the new period value is not available ("NA")
$ws.Range("D12").Value = "NA"
$ws.Range("D29").Value = "NA"
$ws.Range("D44").Value = "NA"
